$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.764.15"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.659.59"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "604.42"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").Value = "148.06"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").Value = "5.61"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "27.63"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "3.136.06"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "63.612.44"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "2.658.67"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "11.52"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("D20").Value = "342.99"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "7.00"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "5.59"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").Value = "66.90"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "1.71"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").Value = "9.17"
$ws.Range("E26").Value = "  +9.44%  "
$ws.Range("D27").Value = "1.54"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "555.58"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.165"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +4.05%  "
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").Value = "0.0₃0822"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  +7.84%  "
$ws.Range("D36").Value = "167.48"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  +9.22%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "19.18"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "169.15"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "3.79"
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D44").Value = "22.63"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "0.0578"
$ws.Range("E45").Value = "  +5.28%  "
$ws.Range("D46").Value = "0.632"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "0.0248"
$ws.Range("E47").Value = "  +4.56%  "
$ws.Range("D48").Value = "0.0965"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "18.93"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "1.87"
$ws.Range("E50").Value = "  +10.86%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.178"
$ws.Range("E51").Value = "  +3.32%  "
